# feat: add 2022-Q3 data
#
# Create the new "2022-Q3" sheet by duplicating the "2022-Q2" sheet (same
# column layout/styling), placing the copy immediately before "2022-Q2",
# then replace its contents with the 2022-Q3 fund data. Finally, insert a
# matching summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# 2022-Q2 had 14 funds (15 rows incl. header); 2022-Q3 only has 4 funds,
# so drop the now-unused rows 6 through 15.
$newSheet.Range("A6:H15").Clear()

# Columns B-G hold text values (fund code/name/numbers-as-text) in the
# source data; force text formatting so values like "009837" or "14.02"
# are not reinterpreted as numbers.
$newSheet.Range("B2:G5").NumberFormat = "@"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "009837"
$newSheet.Range("C2").Value = "华夏磐锐一年定期开放混合A"
$newSheet.Range("D2").Value = "14.02"
$newSheet.Range("E2").Value = "94.15"
$newSheet.Range("F2").Value = "4.54"
$newSheet.Range("G2").Value = "0.6365"
$newSheet.Range("H2").Value = 4

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "005947"
$newSheet.Range("C3").Value = "德邦民裕进取量化精选灵活配置混合A"
$newSheet.Range("D3").Value = "0.46"
$newSheet.Range("E3").Value = "76.66"
$newSheet.Range("F3").Value = "7.06"
$newSheet.Range("G3").Value = "0.0325"
$newSheet.Range("H3").Value = 2

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "009838"
$newSheet.Range("C4").Value = "华夏磐锐一年定期开放混合C"
$newSheet.Range("D4").Value = "0.39"
$newSheet.Range("E4").Value = "94.15"
$newSheet.Range("F4").Value = "4.54"
$newSheet.Range("G4").Value = "0.0177"
$newSheet.Range("H4").Value = 4

# Row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "005948"
$newSheet.Range("C5").Value = "德邦民裕进取量化精选灵活配置混合C"
$newSheet.Range("D5").Value = "0.08"
$newSheet.Range("E5").Value = "76.66"
$newSheet.Range("F5").Value = "7.06"
$newSheet.Range("G5").Value = "0.0056"
$newSheet.Range("H5").Value = 2

# Update the summary sheet ("总计") to insert the new 2022-Q3 row at the top
# of the data (row 2), pushing all the other quarters down by one row and
# renumbering the index column (A). Work bottom-up so we don't clobber data
# before it has been shifted down.
$summary = $wb.Worksheets.Item("总计")

# A new row 8 is needed; copy the formatting from row 7 (A column uses a
# centered/bold style) before filling in the shifted 2020-Q4 values.
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)

# row7 (2020-Q4) -> row8
$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2020-Q4"
$summary.Range("C8").Value = 2
$summary.Range("D8").Value = 0.01

# row6 (2021-Q1) -> row7
$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q1"
$summary.Range("C7").Value = 4
$summary.Range("D7").Value = 0.03

# row5 (2021-Q3) -> row6
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.13

# row4 (2021-Q4) -> row5
$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 0.12

# row3 (2022-Q1) -> row4
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 7
$summary.Range("D4").Value = 0.6899999999999999

# row2 (2022-Q2) -> row3
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 14
$summary.Range("D3").Value = 1.95

# new row2 (2022-Q3)
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.6899999999999999
